$d = $word.ActiveDocument

# --- Locate the first paragraph ("This is a Microsoft word document.") ---
$p1 = $d.Paragraphs(1)

# Sanity-check we found the right paragraph; fall back to a Find if not.
if ($p1.Range.Text.TrimEnd([char]13) -ne "This is a Microsoft word document.") {
    $found = $d.Content
    $found.Find.Execute("This is a Microsoft word document.", $true, $false, $false,
                         $false, $false, $true, 1, $false, "", 0)
    $p1 = $found.Paragraphs(1)
}

# 1) Append two trailing spaces to the existing sentence (same run/formatting).
$p1.Range.Text = "This is a Microsoft word document.  "

# 2) Insert a new, empty paragraph right after this one - this becomes the
#    blank line that now separates the sentence from "It will be treated..."
$p1.Range.InsertParagraphAfter()

# Re-fetch paragraph 1 (still the sentence paragraph); insert the red,
# parenthetical "version" note just before its end-of-paragraph mark.
$p1 = $d.Paragraphs(1)
$insertAt = $p1.Range.End - 1
$branchName = "main"

$r1 = $d.Range($insertAt, $insertAt)
$r1.InsertAfter("(This is a change " + [char]0x2013 + " Version for branch ")
$r1.Font.Color = 192

$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter($branchName)
$r2.Font.Color = 192

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter(")")
$r3.Font.Color = 192
